$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 38463212
$ws.Range("I58").Value = 173.9
$ws.Range("J58").Value = 166673330
$ws.Range("K58").Value = 521.7
$ws.Range("L58").Value = 500019990
$ws.Range("M58").Value = -371.7
$ws.Range("N58").Value = -500020290

$ws.Range("H76").Value = 12032.333
$ws.Range("J76").Value = 17997
$ws.Range("L76").Value = 17997
$ws.Range("N76").Value = -18627

$ws.Range("H79").Value = 12032.333
$ws.Range("J79").Value = 17997
$ws.Range("L79").Value = 17997
$ws.Range("N79").Value = -20181

$ws.Range("H80").Value = 32836.312
$ws.Range("I80").Value = 13569.25
$ws.Range("J80").Value = 52103.375
$ws.Range("K80").Value = 40707.75
$ws.Range("L80").Value = 156310.125
$ws.Range("M80").Value = -39709.75
$ws.Range("N80").Value = -158306.125

$ws.Range("H83").Value = 32836.312
$ws.Range("I83").Value = 13569.25
$ws.Range("J83").Value = 52103.375
$ws.Range("K83").Value = 122123.25
$ws.Range("L83").Value = 468930.375
$ws.Range("M83").Value = -117131.25
$ws.Range("N83").Value = -478914.375

$ws.Range("H99").Value = 948
$ws.Range("I99").Value = 148
$ws.Range("J99").Value = 1348
$ws.Range("K99").Value = 444
$ws.Range("L99").Value = 4044
$ws.Range("M99").Value = 1054
$ws.Range("N99").Value = -7040

$ws.Range("H112").Value = 5079.2666
$ws.Range("J112").Value = 5706.846
$ws.Range("L112").Value = 17120.538
$ws.Range("N112").Value = -19336.538

$ws.Range("H116").Value = 22733664
$ws.Range("I116").Value = 41671316
$ws.Range("J116").Value = 8480
$ws.Range("K116").Value = 41671316
$ws.Range("L116").Value = 8480
$ws.Range("M116").Value = -41667874
$ws.Range("N116").Value = -15364

$ws.Range("H132").Value = 922.5714
$ws.Range("I132").Value = 963.4
$ws.Range("K132").Value = 2890.2
$ws.Range("M132").Value = -360.1999999999998

$ws.Range("H137").Value = 5621.8
$ws.Range("I137").Value = 4553.875
$ws.Range("K137").Value = 13661.625
$ws.Range("M137").Value = -11111.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 451.33334
$ws.Range("I5").Value = 451
$ws.Range("J5").Value = 451.5
$ws.Range("K5").Value = 451
$ws.Range("L5").Value = 451.5
$ws.Range("M5").Value = -339
$ws.Range("N5").Value = -675.5

$ws.Range("H74").Value = 18906.271
$ws.Range("I74").Value = 22447.66
$ws.Range("J74").Value = 5035.8335
$ws.Range("K74").Value = 22447.66
$ws.Range("L74").Value = 5035.8335
$ws.Range("M74").Value = -21573.66
$ws.Range("N74").Value = -6783.8335

$ws.Range("H77").Value = 18906.271
$ws.Range("I77").Value = 22447.66
$ws.Range("J77").Value = 5035.8335
$ws.Range("K77").Value = 112238.3
$ws.Range("L77").Value = 25179.1675
$ws.Range("M77").Value = -107870.3
$ws.Range("N77").Value = -33915.1675

$ws.Range("H105").Value = 68460
$ws.Range("J105").Value = 68460
$ws.Range("L105").Value = 68460
$ws.Range("N105").Value = -75448

$ws.Range("H122").Value = 6937
$ws.Range("I122").Value = 5858.125
$ws.Range("K122").Value = 17574.375
$ws.Range("M122").Value = -15124.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 451.33334
$ws.Range("I4").Value = 451
$ws.Range("J4").Value = 451.5
$ws.Range("K4").Value = 451
$ws.Range("L4").Value = 451.5
$ws.Range("M4").Value = -336
$ws.Range("N4").Value = -681.5

$ws.Range("H128").Value = 4629.25
$ws.Range("I128").Value = 4629.25
$ws.Range("K128").Value = 13887.75
$ws.Range("M128").Value = -11397.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 5156
$ws.Range("I76").Value = 5156
$ws.Range("K76").Value = 5156
$ws.Range("M76").Value = -4841

$ws.Range("H79").Value = 5156
$ws.Range("I79").Value = 5156
$ws.Range("K79").Value = 5156
$ws.Range("M79").Value = -4064

$ws.Range("H86").Value = 6260143.5
$ws.Range("I86").Value = 10427258
$ws.Range("J86").Value = 9473
$ws.Range("K86").Value = 10427258
$ws.Range("L86").Value = 9473
$ws.Range("M86").Value = -10426135
$ws.Range("N86").Value = -11719

$ws.Range("H89").Value = 6260143.5
$ws.Range("I89").Value = 10427258
$ws.Range("J89").Value = 9473
$ws.Range("K89").Value = 52136290
$ws.Range("L89").Value = 47365
$ws.Range("M89").Value = -52130674
$ws.Range("N89").Value = -58597

$ws.Range("H93").Value = 44280
$ws.Range("I93").Value = 39449
$ws.Range("J93").Value = 53942
$ws.Range("K93").Value = 39449
$ws.Range("L93").Value = 53942
$ws.Range("M93").Value = -37577
$ws.Range("N93").Value = -57686

$ws.Range("H99").Value = 4558.375
$ws.Range("I99").Value = 2835.111
$ws.Range("J99").Value = 6774
$ws.Range("K99").Value = 2835.111
$ws.Range("L99").Value = 6774
$ws.Range("M99").Value = -1337.111
$ws.Range("N99").Value = -9770

$ws.Range("H126").Value = 4558.375
$ws.Range("I126").Value = 2835.111
$ws.Range("J126").Value = 6774
$ws.Range("K126").Value = 8505.332999999999
$ws.Range("L126").Value = 20322
$ws.Range("M126").Value = -6035.332999999999
$ws.Range("N126").Value = -25262

$ws.Range("H134").Value = 6498.1714
$ws.Range("I134").Value = 3335.75
$ws.Range("J134").Value = 8148.1304
$ws.Range("K134").Value = 10007.25
$ws.Range("L134").Value = 24444.3912
$ws.Range("M134").Value = -7472.25
$ws.Range("N134").Value = -29514.3912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 183.5
$ws.Range("I8").Value = 183.5
$ws.Range("K8").Value = 550.5
$ws.Range("M8").Value = -411.5

$ws.Range("H38").Value = 55555600
$ws.Range("J38").Value = 125000060
$ws.Range("L38").Value = 375000180
$ws.Range("N38").Value = -375000874

$ws.Range("H87").Value = 250002370
$ws.Range("I87").Value = 250002370
$ws.Range("K87").Value = 750007110
$ws.Range("M87").Value = -750005862

$ws.Range("H90").Value = 250002370
$ws.Range("I90").Value = 250002370
$ws.Range("K90").Value = 2250021330
$ws.Range("M90").Value = -2250015090

$ws.Range("H132").Value = 4412.875
$ws.Range("I132").Value = 1023.6923
$ws.Range("K132").Value = 9213.2307
$ws.Range("M132").Value = -6683.2307

$ws.Range("H141").Value = 7948.2
$ws.Range("I141").Value = 7948.2
$ws.Range("K141").Value = 23844.6
$ws.Range("M141").Value = -18664.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 16.2
$ws.Range("I2").Value = 16.2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 16.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 96.8
$ws.Range("N2").ClearContents()

$ws.Range("H52").Value = 89998.8
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 89998.8
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 89998.8
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -90516.8

$ws.Range("H102").Value = 2828.7646
$ws.Range("I102").Value = 2718.5186
$ws.Range("J102").Value = 3254
$ws.Range("K102").Value = 2718.5186
$ws.Range("L102").Value = 3254
$ws.Range("M102").Value = -1096.5186
$ws.Range("N102").Value = -6498

$ws.Range("H132").Value = 4116.6924
$ws.Range("I132").Value = 1382.5714
$ws.Range("J132").Value = 15600
$ws.Range("K132").Value = 4147.7142
$ws.Range("L132").Value = 46800
$ws.Range("M132").Value = -1617.7142
$ws.Range("N132").Value = -51860

$ws.Range("H139").Value = 66660
$ws.Range("J139").Value = 66660
$ws.Range("L139").Value = 66660
$ws.Range("N139").Value = -76940

$ws.Range("H140").Value = 70184.664
$ws.Range("J140").Value = 70184.664
$ws.Range("L140").Value = 70184.664
$ws.Range("N140").Value = -80544.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6752
$ws.Range("I7").Value = 5004
$ws.Range("K7").Value = 5004
$ws.Range("M7").Value = -4892

$ws.Range("H40").Value = 6133.273
$ws.Range("I40").Value = 2866.75
$ws.Range("K40").Value = 2866.75
$ws.Range("M40").Value = -2730.75

$ws.Range("H82").Value = 834789.4399999999
$ws.Range("I82").Value = 1667433.5
$ws.Range("J82").Value = 2145.3333
$ws.Range("K82").Value = 1667433.5
$ws.Range("L82").Value = 2145.3333
$ws.Range("M82").Value = -1667072.5
$ws.Range("N82").Value = -2867.3333

$ws.Range("H85").Value = 834789.4399999999
$ws.Range("I85").Value = 1667433.5
$ws.Range("J85").Value = 2145.3333
$ws.Range("K85").Value = 1667433.5
$ws.Range("L85").Value = 2145.3333
$ws.Range("M85").Value = -1666185.5
$ws.Range("N85").Value = -4641.3333

$ws.Range("H93").Value = 851.9048
$ws.Range("I93").Value = 811.8333
$ws.Range("J93").Value = 905.3333
$ws.Range("K93").Value = 811.8333
$ws.Range("L93").Value = 905.3333
$ws.Range("M93").Value = 436.1667
$ws.Range("N93").Value = -3401.3333

$ws.Range("H126").Value = 6752
$ws.Range("I126").Value = 5004
$ws.Range("K126").Value = 15012
$ws.Range("M126").Value = -12542

$ws.Range("H136").Value = 12131.546
$ws.Range("I136").Value = 2364.3
$ws.Range("J136").Value = 20270.916
$ws.Range("K136").Value = 7092.900000000001
$ws.Range("L136").Value = 60812.74800000001
$ws.Range("M136").Value = -4542.900000000001
$ws.Range("N136").Value = -65912.74800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8123.4
$ws.Range("I126").Value = 9402
$ws.Range("K126").Value = 28206
$ws.Range("M126").Value = -25736

$ws.Range("H132").Value = 4899.8286
$ws.Range("I132").Value = 4955.68
$ws.Range("J132").Value = 4760.2
$ws.Range("K132").Value = 14867.04
$ws.Range("L132").Value = 14280.6
$ws.Range("M132").Value = -12337.04
$ws.Range("N132").Value = -19340.6

$ws.Range("H136").Value = 19426296
$ws.Range("I136").Value = 41667630
$ws.Range("J136").Value = 362294.6
$ws.Range("K136").Value = 125002890
$ws.Range("L136").Value = 1086883.8
$ws.Range("M136").Value = -125000340
$ws.Range("N136").Value = -1091983.8
